$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 546-547, shifting existing rows 546-662 down to 548-664
$ws.Rows("546:547").Insert()

# Row 546 — new weekly data point (Primera)
$ws.Cells.Item(546, 1).Value = 7
$ws.Cells.Item(546, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(546, 3).Value = "Ñuble"
$ws.Cells.Item(546, 4).Value = 45258
$ws.Cells.Item(546, 5).Value = 16
$ws.Cells.Item(546, 6).Value = 100112023
$ws.Cells.Item(546, 7).Value = "Brócoli"
$ws.Cells.Item(546, 8).Value = "Sin especificar"
$ws.Cells.Item(546, 9).Value = "Primera"
$ws.Cells.Item(546, 10).Value = 500
$ws.Cells.Item(546, 11).Value = 1300
$ws.Cells.Item(546, 12).Value = 1300
$ws.Cells.Item(546, 13).Value = 1300
$ws.Cells.Item(546, 14).Value = "$/unidad"
$ws.Cells.Item(546, 15).Value = "Región del Maule"
$ws.Cells.Item(546, 16).Value = 1300
$ws.Cells.Item(546, 17).Value = 1
$ws.Cells.Item(546, 18).Value = "Hortaliza"

# Row 547 — new weekly data point (Segunda)
$ws.Cells.Item(547, 1).Value = 7
$ws.Cells.Item(547, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(547, 3).Value = "Ñuble"
$ws.Cells.Item(547, 4).Value = 45258
$ws.Cells.Item(547, 5).Value = 16
$ws.Cells.Item(547, 6).Value = 100112023
$ws.Cells.Item(547, 7).Value = "Brócoli"
$ws.Cells.Item(547, 8).Value = "Sin especificar"
$ws.Cells.Item(547, 9).Value = "Segunda"
$ws.Cells.Item(547, 10).Value = 300
$ws.Cells.Item(547, 11).Value = 1000
$ws.Cells.Item(547, 12).Value = 1000
$ws.Cells.Item(547, 13).Value = 1000
$ws.Cells.Item(547, 14).Value = "$/unidad"
$ws.Cells.Item(547, 15).Value = "Región del Maule"
$ws.Cells.Item(547, 16).Value = 1000
$ws.Cells.Item(547, 17).Value = 1
$ws.Cells.Item(547, 18).Value = "Hortaliza"

# Apply the date number format used elsewhere in column D to the two new date cells
$ws.Range("D546:D547").NumberFormat = "YYYY-MM-DD HH:MM:SS"
